$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.662.31"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "'1.894.55"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'238.88"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'0.9979"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.4790"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("D8").Value = "'0.2830"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.06532"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "'1.987.37"
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("D11").Value = "'0.07480"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "'16.64"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'5.093"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "'88.08"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'0.6660"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "'30.589.49"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.30"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9993"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007586"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'229.64"
$ws.Range("E20").Value = "  +4.50%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.120.65"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.294"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'0.9994"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'6.229"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "'18.50"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "'1.954"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "'1.401"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("D30").Value = "'0.09750"
$ws.Range("E30").Value = "  +5.39%  "
$ws.Range("D31").Value = "'4.348"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'4.013"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'0.05054"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'1.228"
$ws.Range("E34").Value = "  +7.24%  "
$ws.Range("D35").Value = "'0.7542"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'2.711"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.01871"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'2.644"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.085"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9102"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'106.25"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'0.4285"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").Value = "'5.778"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'7.389"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "'64.00"
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").Value = "'0.1272"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.982"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.477"
$ws.Range("E49").Value = "  -6.46%  "
$ws.Range("D50").Value = "'33.68"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "'0.05652"
$ws.Range("E51").Value = "  -1.27%  "
